{"js": "// Update the East-Asian font used by the document's paragraph styles from\n// \"DejaVu Sans\" to \"Tahoma\" (Normal / Heading \u2014 the two styles that already\n// carry an explicit <w:rFonts> block), and make the inherited complex-script\n// font (\"DejaVu Sans\") explicit on the styles that previously had an empty\n// <w:rPr/> (List / Caption / Index). This mirrors the vignette docx style\n// refresh described in the commit.\n\nconst styles = context.document.getStyles();\n\nconst normal = styles.getByNameOrNullObject(\"Normal\");\nconst heading = styles.getByNameOrNullObject(\"Heading\");\nconst list = styles.getByNameOrNullObject(\"List\");\nconst caption = styles.getByNameOrNullObject(\"Caption\");\nconst index = styles.getByNameOrNullObject(\"Index\");\n\nawait context.sync();\n\n// Normal + Heading: switch the east-Asian face to \"Tahoma\" while leaving the\n// ascii/hAnsi/complex-script faces untouched.\nif (!normal.isNullObject) {\n  normal.font.nameFarEast = \"Tahoma\";\n}\nif (!heading.isNullObject) {\n  heading.font.nameFarEast = \"Tahoma\";\n}\n\n// List / Caption / Index: these styles had no direct <w:rFonts> at all, so\n// their complex-script font was only ever inherited. Pin it explicitly to\n// the existing effective value (\"DejaVu Sans\") without touching anything\n// else in their (empty) run properties.\nif (!list.isNullObject) {\n  list.font.nameBidirectional = \"DejaVu Sans\";\n}\nif (!caption.isNullObject) {\n  caption.font.nameBidirectional = \"DejaVu Sans\";\n}\nif (!index.isNullObject) {\n  index.font.nameBidirectional = \"DejaVu Sans\";\n}\n\nawait context.sync();\n", "ps1": "# Update the East-Asian font used by the document's paragraph styles from\n# \"DejaVu Sans\" to \"Tahoma\" (Normal / Heading -- the two styles that already\n# carry an explicit rFonts block), and make the inherited complex-script\n# font (\"DejaVu Sans\") explicit on the styles that previously had an empty\n# rPr (List / Caption / Index). This mirrors the vignette docx style\n# refresh described in the commit.\n\n$d = $word.ActiveDocument\n\n$d.Styles(\"Normal\").Font.NameFarEast = \"Tahoma\"\n$d.Styles(\"Heading\").Font.NameFarEast = \"Tahoma\"\n\n$d.Styles(\"List\").Font.NameBi = \"DejaVu Sans\"\n$d.Styles(\"Caption\").Font.NameBi = \"DejaVu Sans\"\n$d.Styles(\"Index\").Font.NameBi = \"DejaVu Sans\"\n"}
